# Apply the row-content rotation described by the diff.
#
# The rows keep their row numbers (A4, A5, ... stay at rows 4, 5, 7, 8, 9)
# but the *content* of each row cyclically moves to a different row:
#   row4 <-> row5              (swap)
#   row7 -> row9 -> row8 -> row7  (3-cycle)
#
# Column C (Valideringsstatus) is always empty in this sheet, so it is
# untouched. Column I (Antal) holds numbers stored as *text* in the
# original file; assigning a plain digit string gets silently re-coerced
# back to a real number, so a quote-prefixed string ("'5") is used to force
# text, then the style is reset to "Normal" so the quote-prefix flag
# doesn't leave a stray cell style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4,1).Value = 131106823
$ws.Cells.Item(4,2).Value = 80348
$ws.Cells.Item(4,4).Value = "NT"
$ws.Cells.Item(4,5).Value = 6458
$ws.Cells.Item(4,6).Value = "Lunglav"
$ws.Cells.Item(4,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(4,8).Value = "(L.) Hoffm."
$ws.Cells.Item(4,9).Value = "'5"
$ws.Cells.Item(4,9).Style = "Normal"
$ws.Cells.Item(4,10).Value = "dm²"
$ws.Cells.Item(4,17).Value = 600504
$ws.Cells.Item(4,18).Value = 6974579
$ws.Cells.Item(4,24).Value = "2025_0354"
$ws.Cells.Item(4,26).Value = "12:44"
$ws.Cells.Item(4,28).Value = "12:44"
$ws.Cells.Item(5,1).Value = 131106812
$ws.Cells.Item(5,2).Value = 99013
$ws.Cells.Item(5,4).Value = "VU"
$ws.Cells.Item(5,5).Value = 220787
$ws.Cells.Item(5,6).Value = "Knärot"
$ws.Cells.Item(5,7).Value = "Goodyera repens"
$ws.Cells.Item(5,8).Value = "(L.) R. Br."
$ws.Cells.Item(5,9).Value = "'20"
$ws.Cells.Item(5,9).Style = "Normal"
$ws.Cells.Item(5,10).Value = "plantor/tuvor"
$ws.Cells.Item(5,17).Value = 600509
$ws.Cells.Item(5,18).Value = 6974463
$ws.Cells.Item(5,24).Value = "2025_0365"
$ws.Cells.Item(5,26).Value = "13:32"
$ws.Cells.Item(5,28).Value = "13:32"
$ws.Cells.Item(7,1).Value = 131106813
$ws.Cells.Item(7,2).Value = 99036
$ws.Cells.Item(7,5).Value = 221952
$ws.Cells.Item(7,6).Value = "Spindelblomster"
$ws.Cells.Item(7,7).Value = "Neottia cordata"
$ws.Cells.Item(7,8).Value = "(L.) Rich."
$ws.Cells.Item(7,9).Value = "'15"
$ws.Cells.Item(7,9).Style = "Normal"
$ws.Cells.Item(7,17).Value = 600498
$ws.Cells.Item(7,18).Value = 6974468
$ws.Cells.Item(7,24).Value = "2025_0364"
$ws.Cells.Item(7,26).Value = "13:31"
$ws.Cells.Item(7,28).Value = "13:31"
$ws.Cells.Item(8,1).Value = 131106821
$ws.Cells.Item(8,2).Value = 80348
$ws.Cells.Item(8,4).Value = "NT"
$ws.Cells.Item(8,5).Value = 6458
$ws.Cells.Item(8,6).Value = "Lunglav"
$ws.Cells.Item(8,7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(8,8).Value = "(L.) Hoffm."
$ws.Cells.Item(8,9).Value = "'1"
$ws.Cells.Item(8,9).Style = "Normal"
$ws.Cells.Item(8,10).Value = "dm²"
$ws.Cells.Item(8,17).Value = 600501
$ws.Cells.Item(8,18).Value = 6974538
$ws.Cells.Item(8,24).Value = "2025_0356"
$ws.Cells.Item(8,26).Value = "12:52"
$ws.Cells.Item(8,28).Value = "12:52"
$ws.Cells.Item(8,29).Value = "På sälglåga"
$ws.Cells.Item(9,1).Value = 131106822
$ws.Cells.Item(9,2).Value = 98930
$ws.Cells.Item(9,4).Value = "LC"
$ws.Cells.Item(9,5).Value = 219790
$ws.Cells.Item(9,6).Value = "Fläcknycklar"
$ws.Cells.Item(9,7).Value = "Dactylorhiza maculata"
$ws.Cells.Item(9,8).Value = "(L.) Soó"
$ws.Cells.Item(9,10).Value = "plantor/tuvor"
$ws.Cells.Item(9,17).Value = 600502
$ws.Cells.Item(9,18).Value = 6974543
$ws.Cells.Item(9,24).Value = "2025_0355"
$ws.Cells.Item(9,26).Value = "12:49"
$ws.Cells.Item(9,28).Value = "12:49"
$ws.Cells.Item(9,29).ClearContents()

